$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Append the new log row (row 44) under the existing data.
$ws.Range("A44").Value = "Wat zijn jullie openingstijden?"
$ws.Range("B44").Value = "mailmind.test@zohomail.eu"
$ws.Range("C44").Value = "Hallo, ik zou graag willen weten wat jullie openingstijden zijn. Dank je wel!"
$ws.Range("D44").Value = "Informatieaanvraag"
$ws.Range("E44").Value = "Beste heer/mevrouw,`nBedankt voor uw bericht. Onze openingstijden zijn maandag tot en met vrijdag van 9.00 uur tot 17.30 uur. Mocht u nog verdere vragen hebben, dan hoor ik het graag.`nMet vriendelijke groet,`n[Naam assistent]"
$ws.Range("F44").Value = "2025-06-17 22:05:03"
$ws.Range("G44").Value = "Ja"

# Extend the conditional formatting ranges to cover the new row.
$ws.Range("D2:D43").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D44"))
$ws.Range("G2:G43").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G44"))

# Update the Dashboard summary count for "Informatieaanvraag".
$dash.Range("B2").Value = 19
